# Auto-generated Excel COM-interop script
# Updates market-price-derived profit columns (H:N) on several rows
# across all 8 job sheets, per the scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item(1)
# Row 18
$ws.Cells.Item(18, 8).Value2 = 3749
$ws.Cells.Item(18, 9).Value2 = 1165.1666
$ws.Cells.Item(18, 10).Value2 = 11500.5
$ws.Cells.Item(18, 11).Value2 = 1165.1666
$ws.Cells.Item(18, 12).Value2 = 11500.5
$ws.Cells.Item(18, 13).Value2 = -881.1666
$ws.Cells.Item(18, 14).Value2 = -12068.5
# Row 74
$ws.Cells.Item(74, 8).Value2 = 7952.7393
$ws.Cells.Item(74, 9).Value2 = 5667.6665
$ws.Cells.Item(74, 10).Value2 = 9421.714
$ws.Cells.Item(74, 11).Value2 = 5667.6665
$ws.Cells.Item(74, 12).Value2 = 9421.714
$ws.Cells.Item(74, 13).Value2 = -4731.6665
$ws.Cells.Item(74, 14).Value2 = -11293.714
# Row 77
$ws.Cells.Item(77, 8).Value2 = 7952.7393
$ws.Cells.Item(77, 9).Value2 = 5667.6665
$ws.Cells.Item(77, 10).Value2 = 9421.714
$ws.Cells.Item(77, 11).Value2 = 28338.3325
$ws.Cells.Item(77, 12).Value2 = 47108.57
$ws.Cells.Item(77, 13).Value2 = -23658.3325
$ws.Cells.Item(77, 14).Value2 = -56468.57
# Row 132
$ws.Cells.Item(132, 8).Value2 = 1717.0566
$ws.Cells.Item(132, 9).Value2 = 1346
$ws.Cells.Item(132, 11).Value2 = 4038
$ws.Cells.Item(132, 13).Value2 = -1508
# Row 133
$ws.Cells.Item(133, 8).Value2 = 79999
$ws.Cells.Item(133, 10).Value2 = 79999
$ws.Cells.Item(133, 12).Value2 = 79999
$ws.Cells.Item(133, 14).Value2 = -90119
# Row 137
$ws.Cells.Item(137, 8).Value2 = 2886.5144
$ws.Cells.Item(137, 9).Value2 = 1909.5834
$ws.Cells.Item(137, 11).Value2 = 5728.7502
$ws.Cells.Item(137, 13).Value2 = -3178.7502
# Row 138
$ws.Cells.Item(138, 8).Value2 = 3680.5178
$ws.Cells.Item(138, 9).Value2 = 1947.5
$ws.Cells.Item(138, 10).Value2 = 4153.159
$ws.Cells.Item(138, 11).Value2 = 5842.5
$ws.Cells.Item(138, 12).Value2 = 12459.477
$ws.Cells.Item(138, 13).Value2 = -702.5
$ws.Cells.Item(138, 14).Value2 = -22739.477

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item(2)
# Row 32
$ws.Cells.Item(32, 8).Value2 = 1225.0779
$ws.Cells.Item(32, 9).Value2 = 1199.9578
$ws.Cells.Item(32, 11).Value2 = 1199.9578
$ws.Cells.Item(32, 13).Value2 = -912.9577999999999
# Row 36
$ws.Cells.Item(36, 8).Value2 = 18578.715
$ws.Cells.Item(36, 9).Value2 = 16683.666
$ws.Cells.Item(36, 11).Value2 = 16683.666
$ws.Cells.Item(36, 13).Value2 = -16337.666
# Row 46
$ws.Cells.Item(46, 8).Value2 = 18332.334
$ws.Cells.Item(46, 9).Value2 = 0
$ws.Cells.Item(46, 11).Value2 = 0
$ws.Cells.Item(46, 13).ClearContents()
# Row 122
$ws.Cells.Item(122, 8).Value2 = 4148.9546
$ws.Cells.Item(122, 9).Value2 = 3406
$ws.Cells.Item(122, 11).Value2 = 10218
$ws.Cells.Item(122, 13).Value2 = -7768
# Row 128
$ws.Cells.Item(128, 8).Value2 = 0
$ws.Cells.Item(128, 10).Value2 = 0
$ws.Cells.Item(128, 12).Value2 = 0
$ws.Cells.Item(128, 14).ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item(3)
# Row 21
$ws.Cells.Item(21, 8).Value2 = 34694.75
$ws.Cells.Item(21, 10).Value2 = 34694.75
$ws.Cells.Item(21, 12).Value2 = 34694.75
$ws.Cells.Item(21, 14).Value2 = -35166.75
# Row 105
$ws.Cells.Item(105, 8).Value2 = 18636.75
$ws.Cells.Item(105, 9).Value2 = 20682.5
$ws.Cells.Item(105, 11).Value2 = 20682.5
$ws.Cells.Item(105, 13).Value2 = -18935.5
# Row 132
$ws.Cells.Item(132, 8).Value2 = 61999.2
$ws.Cells.Item(132, 9).Value2 = 30000
$ws.Cells.Item(132, 10).Value2 = 69999
$ws.Cells.Item(132, 11).Value2 = 30000
$ws.Cells.Item(132, 12).Value2 = 69999
$ws.Cells.Item(132, 13).Value2 = -24940
$ws.Cells.Item(132, 14).Value2 = -80119

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item(4)
# Row 31
$ws.Cells.Item(31, 8).Value2 = 24463.059
$ws.Cells.Item(31, 9).Value2 = 2857.476
$ws.Cells.Item(31, 11).Value2 = 2857.476
$ws.Cells.Item(31, 13).Value2 = -2562.476
# Row 34
$ws.Cells.Item(34, 8).Value2 = 24463.059
$ws.Cells.Item(34, 9).Value2 = 2857.476
$ws.Cells.Item(34, 11).Value2 = 2857.476
$ws.Cells.Item(34, 13).Value2 = -2655.476
# Row 105
$ws.Cells.Item(105, 8).Value2 = 3467.8572
$ws.Cells.Item(105, 9).Value2 = 3390.4443
$ws.Cells.Item(105, 11).Value2 = 3390.4443
$ws.Cells.Item(105, 13).Value2 = -1643.4443
# Row 132
$ws.Cells.Item(132, 8).Value2 = 3532
$ws.Cells.Item(132, 9).Value2 = 2053.9167
$ws.Cells.Item(132, 11).Value2 = 6161.750100000001
$ws.Cells.Item(132, 13).Value2 = -3631.750100000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item(5)
# Row 4
$ws.Cells.Item(4, 8).Value2 = 3275259
$ws.Cells.Item(4, 9).Value2 = 3675437.2
$ws.Cells.Item(4, 10).Value2 = 1594510.5
$ws.Cells.Item(4, 11).Value2 = 11026311.6
$ws.Cells.Item(4, 12).Value2 = 4783531.5
$ws.Cells.Item(4, 13).Value2 = -11026199.6
$ws.Cells.Item(4, 14).Value2 = -4783755.5
# Row 82
$ws.Cells.Item(82, 8).Value2 = 3499.5
$ws.Cells.Item(82, 9).Value2 = 3499.5
$ws.Cells.Item(82, 11).Value2 = 10498.5
$ws.Cells.Item(82, 13).Value2 = -10092.5
# Row 85
$ws.Cells.Item(85, 8).Value2 = 3499.5
$ws.Cells.Item(85, 9).Value2 = 3499.5
$ws.Cells.Item(85, 11).Value2 = 10498.5
$ws.Cells.Item(85, 13).Value2 = -9094.5
# Row 107
$ws.Cells.Item(107, 8).Value2 = 1394.579
$ws.Cells.Item(107, 10).Value2 = 2484.25
$ws.Cells.Item(107, 12).Value2 = 7452.75
$ws.Cells.Item(107, 14).Value2 = -11292.75
# Row 114
$ws.Cells.Item(114, 8).Value2 = 4394.909
$ws.Cells.Item(114, 9).Value2 = 2380.75
$ws.Cells.Item(114, 10).Value2 = 5545.857
$ws.Cells.Item(114, 11).Value2 = 7142.25
$ws.Cells.Item(114, 12).Value2 = 16637.571
$ws.Cells.Item(114, 13).Value2 = -3888.25
$ws.Cells.Item(114, 14).Value2 = -23145.571
# Row 132
$ws.Cells.Item(132, 8).Value2 = 3391.4119
$ws.Cells.Item(132, 10).Value2 = 3922.5334
$ws.Cells.Item(132, 12).Value2 = 35302.8006
$ws.Cells.Item(132, 14).Value2 = -40362.8006
# Row 140
$ws.Cells.Item(140, 8).Value2 = 2086.25
$ws.Cells.Item(140, 9).Value2 = 1276.4117
$ws.Cells.Item(140, 11).Value2 = 3829.2351
$ws.Cells.Item(140, 13).Value2 = 1350.7649

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item(6)
# Row 102
$ws.Cells.Item(102, 8).Value2 = 2153.963
$ws.Cells.Item(102, 9).Value2 = 1528.3125
$ws.Cells.Item(102, 10).Value2 = 3064
$ws.Cells.Item(102, 11).Value2 = 1528.3125
$ws.Cells.Item(102, 12).Value2 = 3064
$ws.Cells.Item(102, 13).Value2 = 93.6875
$ws.Cells.Item(102, 14).Value2 = -6308
# Row 122
$ws.Cells.Item(122, 8).Value2 = 13868.5
$ws.Cells.Item(122, 9).Value2 = 14852.818
$ws.Cells.Item(122, 10).Value2 = 12321.714
$ws.Cells.Item(122, 11).Value2 = 44558.454
$ws.Cells.Item(122, 12).Value2 = 36965.142
$ws.Cells.Item(122, 13).Value2 = -42108.454
$ws.Cells.Item(122, 14).Value2 = -41865.142

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item(7)
# Row 40
$ws.Cells.Item(40, 8).Value2 = 19250.273
$ws.Cells.Item(40, 9).Value2 = 19274.8
$ws.Cells.Item(40, 11).Value2 = 19274.8
$ws.Cells.Item(40, 13).Value2 = -19138.8
# Row 63
$ws.Cells.Item(63, 8).Value2 = 45000
$ws.Cells.Item(63, 9).Value2 = 45000
$ws.Cells.Item(63, 11).Value2 = 45000
$ws.Cells.Item(63, 13).Value2 = -44251
# Row 66
$ws.Cells.Item(66, 8).Value2 = 45000
$ws.Cells.Item(66, 9).Value2 = 45000
$ws.Cells.Item(66, 11).Value2 = 135000
$ws.Cells.Item(66, 13).Value2 = -131256
# Row 122
$ws.Cells.Item(122, 8).Value2 = 7176.4
$ws.Cells.Item(122, 9).Value2 = 5680.5713
$ws.Cells.Item(122, 11).Value2 = 17041.7139
$ws.Cells.Item(122, 13).Value2 = -14591.7139

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item(8)
# Row 22
$ws.Cells.Item(22, 8).Value2 = 7599.4
$ws.Cells.Item(22, 10).Value2 = 6999.25
$ws.Cells.Item(22, 12).Value2 = 6999.25
$ws.Cells.Item(22, 14).Value2 = -7585.25
# Row 132
$ws.Cells.Item(132, 8).Value2 = 4324.4033
$ws.Cells.Item(132, 9).Value2 = 2660.4717
$ws.Cells.Item(132, 10).Value2 = 14123.111
$ws.Cells.Item(132, 11).Value2 = 7981.4151
$ws.Cells.Item(132, 12).Value2 = 42369.333
$ws.Cells.Item(132, 13).Value2 = -5451.4151
$ws.Cells.Item(132, 14).Value2 = -47429.333
# Row 136
$ws.Cells.Item(136, 8).Value2 = 4819.7144
$ws.Cells.Item(136, 9).Value2 = 2075.625
$ws.Cells.Item(136, 11).Value2 = 6226.875
$ws.Cells.Item(136, 13).Value2 = -3676.875

Write-Output "Applied scheduled market-data update to all sheets."
